$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# The workbook currently ends with a single "总计" (totals) sheet.
# We need to:
#   1) Insert a new "2022-Q1" sheet (holding per-fund holdings data) right
#      before the "总计" sheet.
#   2) Update the "总计" sheet so it gets a new first data row for 2022-Q1
#      (and the other rows shift down by one).
# ---------------------------------------------------------------------------

$oldTotal = $wb.Worksheets.Item("总计")

# Reuse an existing, already-styled quarter sheet as a formatting template
# (bold/centered/bordered header row + first column).
$template = $wb.Worksheets.Item("2021-Q4")

# --- Step 1: turn the existing "总计" sheet into the new "2022-Q1" sheet ---
$q1 = $oldTotal
$q1.Cells.Clear()
$q1.Name = "2022-Q1"

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1fundRows = @(
    @("206009", "鹏华新兴产业混合", "44.95", "90.17", "6.79", "3.0521", 2),
    @("010695", "华夏磐益一年定期开放混合", "18.02", "82.41", "2.93", "0.5280", 6),
    @("014125", "华夏中证1000指数增强A", "7.03", "89.75", "0.83", "0.0583", 6),
    @("014126", "华夏中证1000指数增强C", "6.09", "89.75", "0.83", "0.0505", 6)
)

$r = 2
foreach ($row in $q1fundRows) {
    $q1.Range("A$r").Value = $r - 2
    $q1.Range("B$r").Value = "'" + $row[0]
    $q1.Range("C$r").Value = $row[1]
    $q1.Range("D$r").Value = "'" + $row[2]
    $q1.Range("E$r").Value = "'" + $row[3]
    $q1.Range("F$r").Value = "'" + $row[4]
    $q1.Range("G$r").Value = "'" + $row[5]
    $q1.Range("H$r").Value = $row[6]
    $r = $r + 1
}

# Apply the same header / index-column formatting used by the other
# quarterly sheets.
$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$q1.Range("A2:A5").PasteSpecial(-4122)

# --- Step 2: create a brand new "总计" sheet after "2022-Q1" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$total.Name = "总计"

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalRows = @(
    @("2022-Q1", 4, 3.69),
    @("2021-Q4", 4, 4.4),
    @("2021-Q3", 4, 4.73),
    @("2021-Q2", 2, 4.16),
    @("2021-Q1", 4, 4.77),
    @("2020-Q4", 11, 7.8)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Range("A$r").Value = $r - 2
    $total.Range("B$r").Value = $row[0]
    $total.Range("C$r").Value = $row[1]
    $total.Range("D$r").Value = $row[2]
    $r = $r + 1
}

$template.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

$q1.Range("A1").Select() | Out-Null
